# Update LR-pairs_lrc2p/Prok2-Prokr1 sheet with newly computed TPM-based NATMI values.
# A new sending/target cluster, "Resolving-Mac", was added to the analysis, expanding the
# communication table from 3x4 (12 rows) to 3x6 (18 rows) combinations, and every edge
# metric was recalculated against the refreshed TPM expression matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Header row (unchanged) ----
$ws.Range("A1").Value = "Sending cluster"
$ws.Range("B1").Value = "Ligand symbol"
$ws.Range("C1").Value = "Receptor symbol"
$ws.Range("D1").Value = "Target cluster"
$ws.Range("E1").Value = "Ligand-expressing cells"
$ws.Range("F1").Value = "Ligand detection rate"
$ws.Range("G1").Value = "Ligand average expression value"
$ws.Range("H1").Value = "Ligand total expression value"
$ws.Range("I1").Value = "Ligand derived specificity of average expression value"
$ws.Range("J1").Value = "Ligand derived specificity of total expression value"
$ws.Range("K1").Value = "Receptor-expressing cells"
$ws.Range("L1").Value = "Receptor detection rate"
$ws.Range("M1").Value = "Receptor average expression value"
$ws.Range("N1").Value = "Receptor total expression value"
$ws.Range("O1").Value = "Receptor derived specificity of average expression value"
$ws.Range("P1").Value = "Receptor derived specificity of total expression value"
$ws.Range("Q1").Value = "Edge average expression weight"
$ws.Range("R1").Value = "Edge total expression weight"
$ws.Range("S1").Value = "Edge average expression derived specificity"
$ws.Range("T1").Value = "Edge total expression derived specificity"

# Row 2: Inflammatory-Mac -> ECs
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Prok2"
$ws.Range("C2").Value = "Prokr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.130431
$ws.Range("H2").Value = 0.391293
$ws.Range("I2").Value = 0.01061419790058097
$ws.Range("J2").Value = 0.01061419790058097
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.3997403333333334
$ws.Range("N2").Value = 1.199221
$ws.Range("O2").Value = 0.04316820571895996
$ws.Range("P2").Value = 0.04316820571895996
$ws.Range("Q2").Value = 0.052138531417
$ws.Range("R2").Value = 0.469246782753
$ws.Range("S2").Value = 0.0004581958785140323
$ws.Range("T2").Value = 0.0004581958785140323

# Row 3: Inflammatory-Mac -> FAPs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Prok2"
$ws.Range("C3").Value = "Prokr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.130431
$ws.Range("H3").Value = 0.391293
$ws.Range("I3").Value = 0.01061419790058097
$ws.Range("J3").Value = 0.01061419790058097
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 1.29653
$ws.Range("N3").Value = 3.88959
$ws.Range("O3").Value = 0.14001307622399
$ws.Range("P3").Value = 0.14001307622399
$ws.Range("Q3").Value = 0.16910770443
$ws.Range("R3").Value = 1.52196933987
$ws.Range("S3").Value = 0.001486126499710558
$ws.Range("T3").Value = 0.001486126499710558

# Row 4: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Prok2"
$ws.Range("C4").Value = "Prokr1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.130431
$ws.Range("H4").Value = 0.391293
$ws.Range("I4").Value = 0.01061419790058097
$ws.Range("J4").Value = 0.01061419790058097
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 1.446115
$ws.Range("N4").Value = 4.338345
$ws.Range("O4").Value = 0.1561668528484919
$ws.Range("P4").Value = 0.156166852848492
$ws.Range("Q4").Value = 0.188618225565
$ws.Range("R4").Value = 1.697564030085
$ws.Range("S4").Value = 0.001657585881644801
$ws.Range("T4").Value = 0.001657585881644801

# Row 5: Inflammatory-Mac -> MuSCs
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Prok2"
$ws.Range("C5").Value = "Prokr1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1.0
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.130431
$ws.Range("H5").Value = 0.391293
$ws.Range("I5").Value = 0.01061419790058097
$ws.Range("J5").Value = 0.01061419790058097
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.2562553333333333
$ws.Range("N5").Value = 0.7687660000000001
$ws.Range("O5").Value = 0.02767317186552101
$ws.Range("P5").Value = 0.02767317186552102
$ws.Range("Q5").Value = 0.033423639382
$ws.Range("R5").Value = 0.300812754438
$ws.Range("S5").Value = 0.0002937285227174295
$ws.Range("T5").Value = 0.0002937285227174296

# Row 6: Inflammatory-Mac -> Neutrophils
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Prok2"
$ws.Range("C6").Value = "Prokr1"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1.0
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.130431
$ws.Range("H6").Value = 0.391293
$ws.Range("I6").Value = 0.01061419790058097
$ws.Range("J6").Value = 0.01061419790058097
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.7180423333333333
$ws.Range("N6").Value = 2.154127
$ws.Range("O6").Value = 0.07754183547550123
$ws.Range("P6").Value = 0.07754183547550123
$ws.Range("Q6").Value = 0.093654979579
$ws.Range("R6").Value = 0.8428948162109999
$ws.Range("S6").Value = 0.0008230443873112602
$ws.Range("T6").Value = 0.0008230443873112602

# Row 7: Inflammatory-Mac -> Resolving-Mac
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Prok2"
$ws.Range("C7").Value = "Prokr1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1.0
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.130431
$ws.Range("H7").Value = 0.391293
$ws.Range("I7").Value = 0.01061419790058097
$ws.Range("J7").Value = 0.01061419790058097
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 5.143380666666666
$ws.Range("N7").Value = 15.430142
$ws.Range("O7").Value = 0.5554368578675358
$ws.Range("P7").Value = 0.5554368578675359
$ws.Range("Q7").Value = 0.670856283734
$ws.Range("R7").Value = 6.037706553606
$ws.Range("S7").Value = 0.00589551673068289
$ws.Range("T7").Value = 0.005895516730682892

# Row 8: Neutrophils -> ECs
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Prok2"
$ws.Range("C8").Value = "Prokr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 12.078888
$ws.Range("H8").Value = 36.236664
$ws.Range("I8").Value = 0.9829542643309697
$ws.Range("J8").Value = 0.9829542643309698
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 0.3997403333333334
$ws.Range("N8").Value = 1.199221
$ws.Range("O8").Value = 0.04316820571895996
$ws.Range("P8").Value = 0.04316820571895996
$ws.Range("Q8").Value = 4.828418715416
$ws.Range("R8").Value = 43.45576843874401
$ws.Range("S8").Value = 0.04243237189496825
$ws.Range("T8").Value = 0.04243237189496825

# Row 9: Neutrophils -> FAPs
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Prok2"
$ws.Range("C9").Value = "Prokr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 12.078888
$ws.Range("H9").Value = 36.236664
$ws.Range("I9").Value = 0.9829542643309697
$ws.Range("J9").Value = 0.9829542643309698
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 1.29653
$ws.Range("N9").Value = 3.88959
$ws.Range("O9").Value = 0.14001307622399
$ws.Range("P9").Value = 0.14001307622399
$ws.Range("Q9").Value = 15.66064065864
$ws.Range("R9").Value = 140.94576592776
$ws.Range("S9").Value = 0.137626450336468
$ws.Range("T9").Value = 0.1376264503364681

# Row 10: Neutrophils -> Inflammatory-Mac
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Prok2"
$ws.Range("C10").Value = "Prokr1"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 12.078888
$ws.Range("H10").Value = 36.236664
$ws.Range("I10").Value = 0.9829542643309697
$ws.Range("J10").Value = 0.9829542643309698
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 1.446115
$ws.Range("N10").Value = 4.338345
$ws.Range("O10").Value = 0.1561668528484919
$ws.Range("P10").Value = 0.156166852848492
$ws.Range("Q10").Value = 17.46746112012
$ws.Range("R10").Value = 157.20715008108
$ws.Range("S10").Value = 0.1535048739545722
$ws.Range("T10").Value = 0.1535048739545722

# Row 11: Neutrophils -> MuSCs
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "Prok2"
$ws.Range("C11").Value = "Prokr1"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 12.078888
$ws.Range("H11").Value = 36.236664
$ws.Range("I11").Value = 0.9829542643309697
$ws.Range("J11").Value = 0.9829542643309698
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 0.2562553333333333
$ws.Range("N11").Value = 0.7687660000000001
$ws.Range("O11").Value = 0.02767317186552101
$ws.Range("P11").Value = 0.02767317186552102
$ws.Range("Q11").Value = 3.095279470736
$ws.Range("R11").Value = 27.85751523662401
$ws.Range("S11").Value = 0.02720146229277769
$ws.Range("T11").Value = 0.0272014622927777

# Row 12: Neutrophils -> Neutrophils
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "Prok2"
$ws.Range("C12").Value = "Prokr1"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 12.078888
$ws.Range("H12").Value = 36.236664
$ws.Range("I12").Value = 0.9829542643309697
$ws.Range("J12").Value = 0.9829542643309698
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.7180423333333333
$ws.Range("N12").Value = 2.154127
$ws.Range("O12").Value = 0.07754183547550123
$ws.Range("P12").Value = 0.07754183547550123
$ws.Range("Q12").Value = 8.673152923592001
$ws.Range("R12").Value = 78.05837631232801
$ws.Range("S12").Value = 0.0762200778446944
$ws.Range("T12").Value = 0.0762200778446944

# Row 13: Neutrophils -> Resolving-Mac
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "Prok2"
$ws.Range("C13").Value = "Prokr1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 12.078888
$ws.Range("H13").Value = 36.236664
$ws.Range("I13").Value = 0.9829542643309697
$ws.Range("J13").Value = 0.9829542643309698
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 5.143380666666666
$ws.Range("N13").Value = 15.430142
$ws.Range("O13").Value = 0.5554368578675358
$ws.Range("P13").Value = 0.5554368578675359
$ws.Range("Q13").Value = 62.126319014032
$ws.Range("R13").Value = 559.136871126288
$ws.Range("S13").Value = 0.545969028007489
$ws.Range("T13").Value = 0.5459690280074893

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Prok2"
$ws.Range("C14").Value = "Prokr1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1.0
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.079033
$ws.Range("H14").Value = 0.237099
$ws.Range("I14").Value = 0.006431537768449341
$ws.Range("J14").Value = 0.006431537768449341
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 0.3997403333333334
$ws.Range("N14").Value = 1.199221
$ws.Range("O14").Value = 0.04316820571895996
$ws.Range("P14").Value = 0.04316820571895996
$ws.Range("Q14").Value = 0.03159267776433334
$ws.Range("R14").Value = 0.284334099879
$ws.Range("S14").Value = 0.0002776379454776818
$ws.Range("T14").Value = 0.0002776379454776818

# Row 15: Resolving-Mac -> FAPs
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Prok2"
$ws.Range("C15").Value = "Prokr1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1.0
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.079033
$ws.Range("H15").Value = 0.237099
$ws.Range("I15").Value = 0.006431537768449341
$ws.Range("J15").Value = 0.006431537768449341
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 1.29653
$ws.Range("N15").Value = 3.88959
$ws.Range("O15").Value = 0.14001307622399
$ws.Range("P15").Value = 0.14001307622399
$ws.Range("Q15").Value = 0.10246865549
$ws.Range("R15").Value = 0.92221789941
$ws.Range("S15").Value = 0.0009004993878113678
$ws.Range("T15").Value = 0.000900499387811368

# Row 16: Resolving-Mac -> Inflammatory-Mac
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Prok2"
$ws.Range("C16").Value = "Prokr1"
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 1.0
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.079033
$ws.Range("H16").Value = 0.237099
$ws.Range("I16").Value = 0.006431537768449341
$ws.Range("J16").Value = 0.006431537768449341
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 1.446115
$ws.Range("N16").Value = 4.338345
$ws.Range("O16").Value = 0.1561668528484919
$ws.Range("P16").Value = 0.156166852848492
$ws.Range("Q16").Value = 0.114290806795
$ws.Range("R16").Value = 1.028617261155
$ws.Range("S16").Value = 0.001004393012274946
$ws.Range("T16").Value = 0.001004393012274947

# Row 17: Resolving-Mac -> MuSCs
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Prok2"
$ws.Range("C17").Value = "Prokr1"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 1.0
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.079033
$ws.Range("H17").Value = 0.237099
$ws.Range("I17").Value = 0.006431537768449341
$ws.Range("J17").Value = 0.006431537768449341
$ws.Range("K17").Value = 3.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 0.2562553333333333
$ws.Range("N17").Value = 0.7687660000000001
$ws.Range("O17").Value = 0.02767317186552101
$ws.Range("P17").Value = 0.02767317186552102
$ws.Range("Q17").Value = 0.02025262775933333
$ws.Range("R17").Value = 0.182273649834
$ws.Range("S17").Value = 0.0001779810500258881
$ws.Range("T17").Value = 0.0001779810500258881

# Row 18: Resolving-Mac -> Neutrophils
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Prok2"
$ws.Range("C18").Value = "Prokr1"
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 1.0
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.079033
$ws.Range("H18").Value = 0.237099
$ws.Range("I18").Value = 0.006431537768449341
$ws.Range("J18").Value = 0.006431537768449341
$ws.Range("K18").Value = 3.0
$ws.Range("L18").Value = 1.0
$ws.Range("M18").Value = 0.7180423333333333
$ws.Range("N18").Value = 2.154127
$ws.Range("O18").Value = 0.07754183547550123
$ws.Range("P18").Value = 0.07754183547550123
$ws.Range("Q18").Value = 0.05674903973033334
$ws.Range("R18").Value = 0.510741357573
$ws.Range("S18").Value = 0.0004987132434955711
$ws.Range("T18").Value = 0.0004987132434955711

# Row 19: Resolving-Mac -> Resolving-Mac
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Prok2"
$ws.Range("C19").Value = "Prokr1"
$ws.Range("D19").Value = "Resolving-Mac"
$ws.Range("E19").Value = 1.0
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.079033
$ws.Range("H19").Value = 0.237099
$ws.Range("I19").Value = 0.006431537768449341
$ws.Range("J19").Value = 0.006431537768449341
$ws.Range("K19").Value = 3.0
$ws.Range("L19").Value = 1.0
$ws.Range("M19").Value = 5.143380666666666
$ws.Range("N19").Value = 15.430142
$ws.Range("O19").Value = 0.5554368578675358
$ws.Range("P19").Value = 0.5554368578675359
$ws.Range("Q19").Value = 0.4064968042286667
$ws.Range("R19").Value = 3.658471238058
$ws.Range("S19").Value = 0.003572313129363885
$ws.Range("T19").Value = 0.003572313129363886

